# Auto-generated edit script: refresh market-price-derived columns (H-N)
# across multiple Leve-profit worksheets, per scheduled-runner data update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row 137
$ws.Range("H137").Value = 18519870
$ws.Range("I137").Value = 988.7646999999999
$ws.Range("J137").Value = 50001970
$ws.Range("K137").Value = 2966.2941
$ws.Range("L137").Value = 150005910
$ws.Range("M137").Value = -416.2941000000001
$ws.Range("N137").Value = -150011010

$ws = $wb.Worksheets.Item("ARM")
# ARM!row 45
$ws.Range("H45").Value = 1394.0555
$ws.Range("I45").Value = 1258.6923
$ws.Range("J45").Value = 1746
$ws.Range("K45").Value = 1258.6923
$ws.Range("L45").Value = 1746
$ws.Range("M45").Value = -881.6922999999999
$ws.Range("N45").Value = -2500

# ARM!row 61
$ws.Range("H61").Value = 3707203.2
$ws.Range("I61").Value = 4633364.5
$ws.Range("J61").Value = 2558.1667
$ws.Range("K61").Value = 4633364.5
$ws.Range("L61").Value = 2558.1667
$ws.Range("M61").Value = -4633152.5
$ws.Range("N61").Value = -2982.1667

# ARM!row 74
$ws.Range("H74").Value = 23817224
$ws.Range("I74").Value = 41667692
$ws.Range("J74").Value = 16599.334
$ws.Range("K74").Value = 41667692
$ws.Range("L74").Value = 16599.334
$ws.Range("M74").Value = -41666818
$ws.Range("N74").Value = -18347.334

# ARM!row 77
$ws.Range("H77").Value = 23817224
$ws.Range("I77").Value = 41667692
$ws.Range("J77").Value = 16599.334
$ws.Range("K77").Value = 208338460
$ws.Range("L77").Value = 82996.67
$ws.Range("M77").Value = -208334092
$ws.Range("N77").Value = -91732.67

# ARM!row 102
$ws.Range("H102").Value = 1410.6875
$ws.Range("I102").Value = 1321.5385
$ws.Range("J102").Value = 1797
$ws.Range("K102").Value = 1321.5385
$ws.Range("L102").Value = 1797
$ws.Range("M102").Value = 300.4614999999999
$ws.Range("N102").Value = -5041

# ARM!row 136
$ws.Range("H136").Value = 3707203.2
$ws.Range("I136").Value = 4633364.5
$ws.Range("J136").Value = 2558.1667
$ws.Range("K136").Value = 13900093.5
$ws.Range("L136").Value = 7674.500100000001
$ws.Range("M136").Value = -13897543.5
$ws.Range("N136").Value = -12774.5001

$ws = $wb.Worksheets.Item("BSM")
# BSM!row 86
$ws.Range("H86").Value = 1809.5555
$ws.Range("I86").Value = 1897.6666
$ws.Range("J86").Value = 1633.3334
$ws.Range("K86").Value = 1897.6666
$ws.Range("L86").Value = 1633.3334
$ws.Range("M86").Value = -774.6666
$ws.Range("N86").Value = -3879.3334

# BSM!row 89
$ws.Range("H89").Value = 1809.5555
$ws.Range("I89").Value = 1897.6666
$ws.Range("J89").Value = 1633.3334
$ws.Range("K89").Value = 9488.333000000001
$ws.Range("L89").Value = 8166.666999999999
$ws.Range("M89").Value = -3872.333000000001
$ws.Range("N89").Value = -19398.667

# BSM!row 134
$ws.Range("H134").Value = 9144296
$ws.Range("I134").Value = 10587907
$ws.Range("J134").Value = 1426.6666
$ws.Range("K134").Value = 31763721
$ws.Range("L134").Value = 4279.9998
$ws.Range("M134").Value = -31761186
$ws.Range("N134").Value = -9349.9998

$ws = $wb.Worksheets.Item("CRP")
# CRP!row 31
$ws.Range("H31").Value = 2093.5938
$ws.Range("I31").Value = 956.8946999999999
$ws.Range("J31").Value = 3754.923
$ws.Range("K31").Value = 956.8946999999999
$ws.Range("L31").Value = 3754.923
$ws.Range("M31").Value = -661.8946999999999
$ws.Range("N31").Value = -4344.923

# CRP!row 34
$ws.Range("H34").Value = 2093.5938
$ws.Range("I34").Value = 956.8946999999999
$ws.Range("J34").Value = 3754.923
$ws.Range("K34").Value = 956.8946999999999
$ws.Range("L34").Value = 3754.923
$ws.Range("M34").Value = -754.8946999999999
$ws.Range("N34").Value = -4158.923

$ws = $wb.Worksheets.Item("CUL")
# CUL!row 110
$ws.Range("H110").Value = 3600
$ws.Range("I110").Value = 2900
$ws.Range("K110").Value = 8700
$ws.Range("M110").Value = -4610

# CUL!row 113
$ws.Range("H113").Value = 2280.9285
$ws.Range("I113").Value = 503.81482
$ws.Range("J113").Value = 3935.4827
$ws.Range("K113").Value = 1511.44446
$ws.Range("L113").Value = 11806.4481
$ws.Range("M113").Value = 658.5555400000001
$ws.Range("N113").Value = -16146.4481

# CUL!row 131
$ws.Range("H131").Value = 3461.0886
$ws.Range("I131").Value = 6040
$ws.Range("J131").Value = 3129.5144
$ws.Range("K131").Value = 18120
$ws.Range("L131").Value = 9388.5432
$ws.Range("M131").Value = -13080
$ws.Range("N131").Value = -19468.5432

$ws = $wb.Worksheets.Item("GSM")
# GSM!row 132
$ws.Range("H132").Value = 2052.3948
$ws.Range("I132").Value = 1809.4688
$ws.Range("J132").Value = 3348
$ws.Range("K132").Value = 5428.4064
$ws.Range("L132").Value = 10044
$ws.Range("M132").Value = -2898.4064
$ws.Range("N132").Value = -15104

# GSM!row 135
$ws.Range("H135").Value = 38000
$ws.Range("J135").Value = 38000
$ws.Range("L135").Value = 38000
$ws.Range("N135").Value = -48140

$ws = $wb.Worksheets.Item("LTW")
# LTW!row 40
$ws.Range("H40").Value = 1968.5385
$ws.Range("I40").Value = 1882.1666
$ws.Range("J40").Value = 3005
$ws.Range("K40").Value = 1882.1666
$ws.Range("L40").Value = 3005
$ws.Range("M40").Value = -1746.1666
$ws.Range("N40").Value = -3277

# LTW!row 122
$ws.Range("H122").Value = 2308.889
$ws.Range("I122").Value = 1932.3334
$ws.Range("J122").Value = 2416.476
$ws.Range("K122").Value = 5797.0002
$ws.Range("L122").Value = 7249.428
$ws.Range("M122").Value = -3347.0002
$ws.Range("N122").Value = -12149.428

# LTW!row 132
$ws.Range("H132").Value = 4991.7812
$ws.Range("I132").Value = 5242.185
$ws.Range("J132").Value = 3639.6
$ws.Range("K132").Value = 15726.555
$ws.Range("L132").Value = 10918.8
$ws.Range("M132").Value = -13196.555
$ws.Range("N132").Value = -15978.8

# LTW!row 136
$ws.Range("H136").Value = 3467.1428
$ws.Range("I136").Value = 2094.5454
$ws.Range("J136").Value = 8500
$ws.Range("K136").Value = 6283.6362
$ws.Range("L136").Value = 25500
$ws.Range("M136").Value = -3733.6362
$ws.Range("N136").Value = -30600

$ws = $wb.Worksheets.Item("WVR")
# WVR!row 122
$ws.Range("H122").Value = 1265
$ws.Range("I122").Value = 1265
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3795
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1345
$ws.Range("N122").ClearContents()

# WVR!row 132
$ws.Range("H132").Value = 5394.077
$ws.Range("I132").Value = 6232.2
$ws.Range("J132").Value = 2600.3333
$ws.Range("K132").Value = 18696.6
$ws.Range("L132").Value = 7800.999899999999
$ws.Range("M132").Value = -16166.6
$ws.Range("N132").Value = -12860.9999

# WVR!row 136
$ws.Range("H136").Value = 7650.027
$ws.Range("I136").Value = 10234.462
$ws.Range("J136").Value = 1541.3636
$ws.Range("K136").Value = 30703.386
$ws.Range("L136").Value = 4624.0908
$ws.Range("M136").Value = -28153.386
$ws.Range("N136").Value = -9724.0908

